# Refresh the cryptocurrency price/volume table (cols B-E) to match the
# latest scrape, per the commit "Updated cryptos list on Mon Sep 18
# 03:45:16 UTC 2023 with GitHub Actions". Ranks in column A are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.800.58'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '1.638.68'
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("D5").Value = '''218.64'
$ws.Range("E5").Value = '  +0.86%  '

$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D10").Value = '''19.29'
$ws.Range("E10").Value = '  +0.52%  '

$ws.Range("D11").Value = '''0.0844'
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").Value = '1.867.84'
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").Value = '1.635.01'
$ws.Range("E13").Value = '  -0.54%  '

$ws.Range("E14").Value = '  -0.69%  '

$ws.Range("D15").Value = '''0.526'
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").Value = '''64.85'
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("D17").Value = '26.797.80'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").Value = '0.0₃0734'
$ws.Range("E18").Value = '  -0.63%  '

$ws.Range("D19").Value = '''216.08'
$ws.Range("E19").Value = '  +0.83%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").Value = '''6.55'
$ws.Range("E22").Value = '  +4.18%  '

$ws.Range("D23").Value = '''2.36'
$ws.Range("E23").Value = '  -2.78%  '

$ws.Range("E24").Value = '  -2.19%  '

$ws.Range("D25").Value = '''147.24'
$ws.Range("E25").Value = '  +1.61%  '

$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("E28").Value = '  -0.94%  '

$ws.Range("D29").Value = '''15.72'
$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("D30").Value = '''0.0505'
$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("E31").Value = '  +1.39%  '

$ws.Range("D32").Value = '''3.38'
$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("D33").Value = '''2.98'
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("D35").Value = '1.262.02'
$ws.Range("E35").Value = '  -2.27%  '

$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("E38").Value = '  -2.18%  '

$ws.Range("D39").Value = '''0.815'
$ws.Range("E39").Value = '  -1.31%  '

$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").Value = '''0.805'
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("E42").Value = '  -0.29%  '

$ws.Range("D43").Value = '1.778.69'
$ws.Range("E43").Value = '  -0.79%  '

$ws.Range("D44").Value = '''2.13'
$ws.Range("E44").Value = '  -4.54%  '

$ws.Range("E45").Value = '  +0.64%  '

$ws.Range("D46").Value = '''60.24'
$ws.Range("E46").Value = '  +0.47%  '

$ws.Range("E47").Value = '  -2.34%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0515'
$ws.Range("E49").Value = '  -1.05%  '

$ws.Range("D50").Value = '''7.55'
$ws.Range("E50").Value = '  -1.94%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.0962'
$ws.Range("E51").Value = '  -1.49%  '
